# Update results data (row 2) on sheets "2025", "2030", "2035"
# with freshly computed values from the server, per commit "ADD results from server".

$wb = $excel.ActiveWorkbook

# Sheet "2025" (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2900.628494009765
$ws.Range("E2").Value = 290490.7128553876
$ws.Range("G2").Value = 80959.25712662016
$ws.Range("I2").Value = 149670.3797976358
$ws.Range("L2").Value = 509990.4857836801
$ws.Range("M2").Value = 112287.0813999
$ws.Range("N2").Value = 71616.34392528524
$ws.Range("O2").Value = 66590.6101513461

# Sheet "2030" (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 1146.922346758181
$ws.Range("B2").Value = 35028.06713949212
$ws.Range("E2").Value = 164950.6135955845
$ws.Range("I2").Value = 160289.5552232918
$ws.Range("L2").Value = 97345.83303624866
$ws.Range("M2").Value = 61409.586608832
$ws.Range("N2").Value = 20784.29964156081
$ws.Range("O2").Value = 11636.11273122574

# Sheet "2035" (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 20355.53802006342
$ws.Range("B2").Value = 19095.16039105755
$ws.Range("E2").Value = 120635.3076705246
$ws.Range("I2").Value = 169865.8582240109
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 58893.31251480614
$ws.Range("N2").Value = 43001.80408654805
$ws.Range("O2").Value = 52185.94545768837
